$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the value currently in A2 (the "zombie" being moved to the end)
$movedValue = $ws.Range("A2").Value2

# Remove row 2, shifting all rows below it up by one
$ws.Rows(2).Delete() | Out-Null

# Place the moved value into the new last data row (A56)
$ws.Range("A56").Value = $movedValue
